$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.463.43"
$ws.Range("E2").Value = "  +0.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.969.47"
$ws.Range("E3").Value = "  +3.70%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.98"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4655"
$ws.Range("E7").Value = "  +0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3911"
$ws.Range("E8").Value = "  -0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.34"
$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07922"
$ws.Range("E10").Value = "  +0.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9880"
$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.75"
$ws.Range("E12").Value = "  +4.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.962.00"
$ws.Range("E13").Value = "  +3.02%  "

$ws.Range("E14").Value = "  +1.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.827"
$ws.Range("E15").Value = "  +1.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07082"
$ws.Range("E16").Value = "  +1.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.70"
$ws.Range("E17").Value = "  -0.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009933"
$ws.Range("E19").Value = "  -0.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.25"

$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.466.07"
$ws.Range("E22").Value = "  +0.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.529"
$ws.Range("E23").Value = "  +4.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.212.35"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.106"
$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.20"
$ws.Range("E27").Value = "  +1.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.51"
$ws.Range("E28").Value = "  +0.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.776"
$ws.Range("E29").Value = "  -4.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.59"
$ws.Range("E30").Value = "  +0.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.902"
$ws.Range("E31").Value = "  +0.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09394"
$ws.Range("E32").Value = "  +0.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8906"
$ws.Range("E33").Value = "  -1.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.228"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.316"
$ws.Range("E35").Value = "  -0.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.165"
$ws.Range("E36").Value = "  -1.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05812"
$ws.Range("E37").Value = "  +0.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.169"
$ws.Range("E38").Value = "  -1.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02098"
$ws.Range("E39").Value = "  +0.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.740"
$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5710"
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000003089"
$ws.Range("E42").Value = "  +46.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1793"
$ws.Range("E43").Value = "  +0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.639"
$ws.Range("E44").Value = "  -0.75%  "

$ws.Range("E45").Value = "  +7.48%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.78"
$ws.Range("E46").Value = "  -0.89%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5337"
$ws.Range("E47").Value = "  -0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.194"
$ws.Range("E48").Value = "  +0.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06913"
$ws.Range("E49").Value = "  -1.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.827"
$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.56"
$ws.Range("E51").Value = "  +0.47%  "
